{"js": "// Update the 25 two-digit-division answers in the worksheet table.\n// The table has 20 rows x 5 columns; only every 4th row (0, 4, 8, 12, 16)\n// holds the \"a\u00f7b=c, d\" answer text, the others are blank spacer rows.\n// Each populated cell's text is replaced with the new value below, in the\n// same row-major (row then column) order the cells appear in the document.\n\nconst replacements = [\n  \"91\u00f77=13, 0\", \"92\u00f73=30, 2\",\n  \"32\u00f78=4, 0\", \"92\u00f72=46, 0\",\n  \"72\u00f79=8, 0\", \"64\u00f75=12, 4\",\n  \"41\u00f73=13, 2\", \"76\u00f74=19, 0\",\n  \"77\u00f75=15, 2\", \"45\u00f72=22, 1\",\n  \"20\u00f75=4, 0\", \"31\u00f77=4, 3\",\n  \"86\u00f72=43, 0\", \"35\u00f78=4, 3\",\n  \"68\u00f77=9, 5\", \"93\u00f75=18, 3\",\n  \"95\u00f78=11, 7\", \"83\u00f77=11, 6\",\n  \"77\u00f75=15, 2\", \"25\u00f75=5, 0\",\n  \"56\u00f74=14, 0\", \"58\u00f74=14, 2\",\n  \"98\u00f79=10, 8\", \"65\u00f72=32, 1\",\n  \"88\u00f78=11, 0\", \"92\u00f78=11, 4\",\n  \"55\u00f73=18, 1\", \"47\u00f75=9, 2\",\n  \"19\u00f77=2, 5\", \"83\u00f79=9, 2\",\n  \"51\u00f73=17, 0\", \"65\u00f73=21, 2\",\n  \"60\u00f77=8, 4\", \"73\u00f76=12, 1\",\n  \"17\u00f76=2, 5\", \"19\u00f79=2, 1\",\n  \"41\u00f75=8, 1\", \"25\u00f73=8, 1\",\n  \"29\u00f72=14, 1\", \"21\u00f79=2, 3\",\n  \"74\u00f74=18, 2\", \"47\u00f79=5, 2\",\n  \"56\u00f74=14, 0\", \"87\u00f79=9, 6\",\n  \"75\u00f76=12, 3\", \"72\u00f75=14, 2\",\n  \"28\u00f77=4, 0\", \"60\u00f73=20, 0\",\n  \"16\u00f73=5, 1\", \"15\u00f74=3, 3\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every cell's current text so we can walk the table in document order\n// and verify we are rewriting the expected old value before swapping it.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.load(\"value\");\n  }\n}\nawait context.sync();\n\nlet pos = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const text = cell.value;\n    if (text === \"\" || text === undefined) {\n      continue; // blank spacer cell, nothing to update\n    }\n    if (pos >= replacements.length) {\n      continue;\n    }\n    const oldVal = replacements[pos];\n    const newVal = replacements[pos + 1];\n    pos += 2;\n    if (text === oldVal) {\n      cell.value = newVal;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Update the 25 two-digit-division answers in the worksheet table.\n# The table has 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17 in\n# 1-based Word COM indexing) holds the \"a\u00f7b=c, d\" answer text, the rest are\n# blank spacer rows. Each populated cell's text is replaced with the new\n# value below, walking the cells in the same row-major (row then column)\n# order they appear in the document.\n\n$replacements = @(\n  \"91\u00f77=13, 0\", \"92\u00f73=30, 2\",\n  \"32\u00f78=4, 0\", \"92\u00f72=46, 0\",\n  \"72\u00f79=8, 0\", \"64\u00f75=12, 4\",\n  \"41\u00f73=13, 2\", \"76\u00f74=19, 0\",\n  \"77\u00f75=15, 2\", \"45\u00f72=22, 1\",\n  \"20\u00f75=4, 0\", \"31\u00f77=4, 3\",\n  \"86\u00f72=43, 0\", \"35\u00f78=4, 3\",\n  \"68\u00f77=9, 5\", \"93\u00f75=18, 3\",\n  \"95\u00f78=11, 7\", \"83\u00f77=11, 6\",\n  \"77\u00f75=15, 2\", \"25\u00f75=5, 0\",\n  \"56\u00f74=14, 0\", \"58\u00f74=14, 2\",\n  \"98\u00f79=10, 8\", \"65\u00f72=32, 1\",\n  \"88\u00f78=11, 0\", \"92\u00f78=11, 4\",\n  \"55\u00f73=18, 1\", \"47\u00f75=9, 2\",\n  \"19\u00f77=2, 5\", \"83\u00f79=9, 2\",\n  \"51\u00f73=17, 0\", \"65\u00f73=21, 2\",\n  \"60\u00f77=8, 4\", \"73\u00f76=12, 1\",\n  \"17\u00f76=2, 5\", \"19\u00f79=2, 1\",\n  \"41\u00f75=8, 1\", \"25\u00f73=8, 1\",\n  \"29\u00f72=14, 1\", \"21\u00f79=2, 3\",\n  \"74\u00f74=18, 2\", \"47\u00f79=5, 2\",\n  \"56\u00f74=14, 0\", \"87\u00f79=9, 6\",\n  \"75\u00f76=12, 3\", \"72\u00f75=14, 2\",\n  \"28\u00f77=4, 0\", \"60\u00f73=20, 0\",\n  \"16\u00f73=5, 1\", \"15\u00f74=3, 3\"\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$pos = 0\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n  for ($c = 1; $c -le $table.Columns.Count; $c++) {\n    $cell = $table.Cell($r, $c)\n    # Cell text includes a trailing end-of-cell marker (chr 13 + chr 7);\n    # strip it before comparing/replacing.\n    $raw = $cell.Range.Text\n    $text = $raw.TrimEnd([char]7).TrimEnd([char]13)\n    if ([string]::IsNullOrEmpty($text)) {\n      continue\n    }\n    if ($pos -ge $replacements.Length) {\n      continue\n    }\n    $oldVal = $replacements[$pos]\n    $newVal = $replacements[$pos + 1]\n    $pos += 2\n    if ($text -eq $oldVal) {\n      $cell.Range.Text = $newVal\n    }\n  }\n}\n"}
